$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'40.114.73"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = "'  +0.12%  "
$ws.Range('E2').ClearFormats()
$ws.Range('D3').Value = "'2.224.01"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = "'  -0.73%  "
$ws.Range('E3').ClearFormats()
$ws.Range('E4').Value = "'  +0.06%  "
$ws.Range('E4').ClearFormats()
$ws.Range('D5').Value = "'290.66"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = "'  -1.17%  "
$ws.Range('E5').ClearFormats()
$ws.Range('D6').Value = "'88.14"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = "'  +1.79%  "
$ws.Range('E6').ClearFormats()
$ws.Range('D7').Value = "'0.513"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = "'  -0.48%  "
$ws.Range('E7').ClearFormats()
$ws.Range('E8').Value = "'  -0.08%  "
$ws.Range('E8').ClearFormats()
$ws.Range('E9').Value = "'  +0.46%  "
$ws.Range('E9').ClearFormats()
$ws.Range('D10').Value = "'30.55"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = "'  -0.23%  "
$ws.Range('E10').ClearFormats()
$ws.Range('E11').Value = "'  -2.06%  "
$ws.Range('E11').ClearFormats()
$ws.Range('D12').Value = "'0.110"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = "'  +2.91%  "
$ws.Range('E12').ClearFormats()
$ws.Range('D13').Value = "'6.50"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = "'  +1.56%  "
$ws.Range('E13').ClearFormats()
$ws.Range('D14').Value = "'2.568.30"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = "'  -0.79%  "
$ws.Range('E14').ClearFormats()
$ws.Range('D15').Value = "'14.00"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = "'  -1.76%  "
$ws.Range('E15').ClearFormats()
$ws.Range('D16').Value = "'2.222.59"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = "'  -0.53%  "
$ws.Range('E16').ClearFormats()
$ws.Range('E17').Value = "'  +0.59%  "
$ws.Range('E17').ClearFormats()
$ws.Range('D18').Value = "'40.057.07"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = "'  +0.16%  "
$ws.Range('E18').ClearFormats()
$ws.Range('D19').Value = "'11.53"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = "'  +7.35%  "
$ws.Range('E19').ClearFormats()
$ws.Range('E20').Value = "'  -0.93%  "
$ws.Range('E20').ClearFormats()
$ws.Range('D21').Value = "'5.83"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = "'  +0.29%  "
$ws.Range('E21').ClearFormats()
$ws.Range('D22').Value = "'65.76"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = "'  +0.20%  "
$ws.Range('E22').ClearFormats()
$ws.Range('D23').Value = "'236.21"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = "'  +0.69%  "
$ws.Range('E23').ClearFormats()
$ws.Range('E24').Value = "'  +0.00%  "
$ws.Range('E24').ClearFormats()
$ws.Range('D25').Value = "'2.47"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = "'  +1.51%  "
$ws.Range('E25').ClearFormats()
$ws.Range('E26').Value = "'  -0.95%  "
$ws.Range('E26').ClearFormats()
$ws.Range('D27').Value = "'22.67"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = "'  -1.88%  "
$ws.Range('E27').ClearFormats()
$ws.Range('D28').Value = "'2.10"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = "'  -4.59%  "
$ws.Range('E28').ClearFormats()
$ws.Range('E29').Value = "'  -0.30%  "
$ws.Range('E29').ClearFormats()
$ws.Range('D30').Value = "'155.56"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = "'  +0.13%  "
$ws.Range('E30').ClearFormats()
$ws.Range('D31').Value = "'31.91"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = "'  -5.02%  "
$ws.Range('E31').ClearFormats()
$ws.Range('E32').Value = "'  -0.10%  "
$ws.Range('E32').ClearFormats()
$ws.Range('D33').Value = "'4.95"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = "'  +1.78%  "
$ws.Range('E33').ClearFormats()
$ws.Range('D34').Value = "'0.0719"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = "'  +1.02%  "
$ws.Range('E34').ClearFormats()
$ws.Range('D35').Value = "'2.36"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = "'  -0.67%  "
$ws.Range('E35').ClearFormats()
$ws.Range('E36').Value = "'  +6.33%  "
$ws.Range('E36').ClearFormats()
$ws.Range('E37').Value = "'  -0.38%  "
$ws.Range('E37').ClearFormats()
$ws.Range('D38').Value = "'15.81"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = "'  -4.76%  "
$ws.Range('E38').ClearFormats()
$ws.Range('E39').Value = "'  -0.14%  "
$ws.Range('E39').ClearFormats()
$ws.Range('D40').Value = "'1.70"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = "'  +1.58%  "
$ws.Range('E40').ClearFormats()
$ws.Range('D41').Value = "'2.117.18"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = "'  +8.13%  "
$ws.Range('E41').ClearFormats()
$ws.Range('D42').Value = "'3.85"
$ws.Range('D42').ClearFormats()
$ws.Range('E43').Value = "'  -2.16%  "
$ws.Range('E43').ClearFormats()
$ws.Range('B44').Value = "'FraxShare"
$ws.Range('B44').ClearFormats()
$ws.Range('C44').Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range('C44').ClearFormats()
$ws.Range('D44').Value = "'10.00"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = "'  +4.83%  "
$ws.Range('E44').ClearFormats()
$ws.Range('B45').Value = "'VeChain"
$ws.Range('B45').ClearFormats()
$ws.Range('C45').Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range('C45').ClearFormats()
$ws.Range('D45').Value = "'0.0268"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = "'  -1.17%  "
$ws.Range('E45').ClearFormats()
$ws.Range('D46').Value = "'17.85"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = "'  +9.51%  "
$ws.Range('E46').ClearFormats()
$ws.Range('D47').Value = "'2.67"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = "'  +1.47%  "
$ws.Range('E47').ClearFormats()
$ws.Range('D48').Value = "'2.433.24"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = "'  -1.01%  "
$ws.Range('E48').ClearFormats()
$ws.Range('D49').Value = "'89.06"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = "'  -0.23%  "
$ws.Range('E49').ClearFormats()
$ws.Range('E50').Value = "'  -1.50%  "
$ws.Range('E50').ClearFormats()
$ws.Range('E51').Value = "'  -2.78%  "
$ws.Range('E51').ClearFormats()

Write-Output "applied changes"
